$wb = $excel.ActiveWorkbook

# --- Sheet 1: covid19_cases_switzerland ---
$ws1 = $wb.Worksheets.Item("covid19_cases_switzerland")

# New data points added for GE (col I), and BS/JU (cols G/L)
$ws1.Range("I11").Value = 294
$ws1.Range("I12").Value = 382
$ws1.Range("G13").Value = 165
$ws1.Range("L13").Value = 23

# --- Sheet 2: Quellen (sources) ---
$ws2 = $wb.Worksheets.Item("Quellen")

# Row 12: BS
$ws2.Range("A12").Value = "BS"
$ws2.Range("B12").Value = "https://www.coronavirus.bs.ch"
$ws2.Hyperlinks.Add($ws2.Range("B12"), "https://www.coronavirus.bs.ch")
$ws2.Range("B12").Style = "Hyperlink"

# Row 13: JU
$ws2.Range("A13").Value = "JU"
$ws2.Range("B13").Value = "https://www.jura.ch/fr/Autorites/Coronavirus/Accueil/Coronavirus-Informations-officielles-a-la-population-jurassienne.html"
$ws2.Hyperlinks.Add($ws2.Range("B13"), "https://www.jura.ch/fr/Autorites/Coronavirus/Accueil/Coronavirus-Informations-officielles-a-la-population-jurassienne.html")
$ws2.Range("B13").Style = "Hyperlink"

# Row 14: GE
$ws2.Range("A14").Value = "GE"
$ws2.Range("B14").Value = "https://www.ge.ch/document/point-coronavirus-maladie-covid-19"
$ws2.Hyperlinks.Add($ws2.Range("B14"), "https://www.ge.ch/document/point-coronavirus-maladie-covid-19")
$ws2.Range("B14").Style = "Hyperlink"

# --- Selections (cursor position saved with the workbook) ---
$ws2.Activate()
$ws2.Range("C20").Select()

$ws1.Activate()
$ws1.Range("W22").Select()

# --- Window position (best-effort; some hosts may not persist this) ---
$aw = $excel.ActiveWindow
$aw.Left = 6555
$aw.Top = 3045
